# Fruta / hortaliza, semanal
# Update the weekly Frutilla price rows for "Vega Monumental Concepción":
#  - shift the existing weekly buckets (rows 404-411) one slot later,
#  - insert a brand-new latest week (dated 44890) at rows 401-403,
#  - and append the freed week (dated 44496) as new rows 412-414.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (date serial, quality, volume, min, max, avg, price/kg)
$rows = @{
    401 = @(44890, "Especial", 200, 9000,  9000,  9000,  1286)
    402 = @(44890, "Primera",  200, 7500,  7500,  7500,  1071)
    403 = @(44890, "Segunda",  100, 6000,  6000,  6000,  857)
    404 = @(44225, "Especial", 100, 10000, 10000, 10000, 1429)
    405 = @(44225, "Primera",  100, 8000,  8000,  8000,  1143)
    406 = @(44225, "Segunda",  50,  6000,  6000,  6000,  857)
    407 = @(44236, "Especial", 100, 10000, 10000, 10000, 1429)
    408 = @(44236, "Primera",  100, 8000,  8000,  8000,  1143)
    409 = @(44236, "Segunda",  100, 7000,  7000,  7000,  1000)
    410 = @(44335, "Especial", 100, 14000, 14000, 14000, 2000)
    411 = @(44335, "Primera",  100, 12000, 12000, 12000, 1714)
    412 = @(44496, "Especial", 100, 10000, 10000, 10000, 1429)
    413 = @(44496, "Primera",  100, 8000,  8000,  8000,  1143)
    414 = @(44496, "Segunda",  50,  6000,  6000,  6000,  857)
}

foreach ($r in 401..414) {
    $data = $rows[$r]

    # Columns A, B, C, E-K, Q, R, T are identical across all these rows;
    # fill them in for the brand-new rows (412-414) and leave the rest as-is.
    $ws.Cells.Item($r, 1).Value = 11
    $ws.Cells.Item($r, 2).Value = "Vega Monumental Concepción"
    $ws.Cells.Item($r, 3).Value = "Bíobío"
    $ws.Cells.Item($r, 4).Value = $data[0]
    $ws.Cells.Item($r, 5).Value = 8
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100101
    $ws.Cells.Item($r, 8).Value = "Berries"
    $ws.Cells.Item($r, 9).Value = 100112025
    $ws.Cells.Item($r, 10).Value = "Frutilla"
    $ws.Cells.Item($r, 11).Value = "Sin especificar"
    $ws.Cells.Item($r, 12).Value = $data[1]
    $ws.Cells.Item($r, 13).Value = $data[2]
    $ws.Cells.Item($r, 14).Value = $data[3]
    $ws.Cells.Item($r, 15).Value = $data[4]
    $ws.Cells.Item($r, 16).Value = $data[5]
    $ws.Cells.Item($r, 17).Value = "`$/bandeja 7 kilos"
    $ws.Cells.Item($r, 18).Value = "Provincia de Melipilla"
    $ws.Cells.Item($r, 19).Value = $data[6]
    $ws.Cells.Item($r, 20).Value = 7

    # Rows 412-414 are brand new; give column D the same date number format
    # that the rest of the column already carries (rows <= 411 keep theirs).
    if ($r -ge 412) {
        $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    }
}
